$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "update by yangguang2"
$ws.Range("A10").Select()
